# Applies the commit "Atualizado por script em 02-12-2023 20:45":
#  - swap the F:V (match detail) contents between four pairs of rows
#    whose A:E (index/league/date) columns stayed untouched
#  - append three brand-new match rows (122, 123, 124) at the bottom,
#    growing the used range from A1:V121 to A1:V124

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowDetails($r1, $r2) {
    $rng1 = $ws.Range("F$r1`:V$r1")
    $rng2 = $ws.Range("F$r2`:V$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Swap-RowDetails 2 3
Swap-RowDetails 38 39
Swap-RowDetails 77 78
Swap-RowDetails 84 85

# --- Append the three new rows, copying formatting from the last
#     existing row (121) and then filling in the new values ---

$ws.Range("A121:V121").Copy()
$ws.Range("A122:V124").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A122").Value2 = 121
$ws.Range("B122").Value2 = "serbia"
$ws.Range("C122").Value2 = "super-liga"
$ws.Range("D122").Value2 = "2023-2024"
$ws.Range("E122").Value2 = 45262.64583333334
$ws.Range("F122").Value2 = "Mladost"
$ws.Range("G122").Value2 = 1
$ws.Range("H122").Value2 = "Zeleznicar Pancevo"
$ws.Range("I122").Value2 = 0
$ws.Range("J122").Value2 = 2.32
$ws.Range("K122").Value2 = "01/12/2023 03:43"
$ws.Range("L122").Value2 = 2.14
$ws.Range("M122").Value2 = "02/12/2023 15:03"
$ws.Range("N122").Value2 = 2.97
$ws.Range("O122").Value2 = "01/12/2023 03:43"
$ws.Range("P122").Value2 = 3.24
$ws.Range("Q122").Value2 = "02/12/2023 15:03"
$ws.Range("R122").Value2 = 2.94
$ws.Range("S122").Value2 = "01/12/2023 03:43"
$ws.Range("T122").Value2 = 3.38
$ws.Range("U122").Value2 = "02/12/2023 15:03"
$ws.Range("V122").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/mladost-lucani-zeleznicar-pancevo/zka4Yd9c/"

$ws.Range("A123").Value2 = 122
$ws.Range("B123").Value2 = "serbia"
$ws.Range("C123").Value2 = "super-liga"
$ws.Range("D123").Value2 = "2023-2024"
$ws.Range("E123").Value2 = 45262.64583333334
$ws.Range("F123").Value2 = "Partizan"
$ws.Range("G123").Value2 = 3
$ws.Range("H123").Value2 = "Vojvodina"
$ws.Range("I123").Value2 = 1
$ws.Range("J123").Value2 = 1.5
$ws.Range("K123").Value2 = "01/12/2023 03:43"
$ws.Range("L123").Value2 = 1.51
$ws.Range("M123").Value2 = "02/12/2023 15:26"
$ws.Range("N123").Value2 = 4.01
$ws.Range("O123").Value2 = "01/12/2023 03:43"
$ws.Range("P123").Value2 = 3.71
$ws.Range("Q123").Value2 = "02/12/2023 15:26"
$ws.Range("R123").Value2 = 5.3
$ws.Range("S123").Value2 = "01/12/2023 03:43"
$ws.Range("T123").Value2 = 7.24
$ws.Range("U123").Value2 = "02/12/2023 15:26"
$ws.Range("V123").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/partizan-vojvodina/tUjxVN25/"

$ws.Range("A124").Value2 = 123
$ws.Range("B124").Value2 = "serbia"
$ws.Range("C124").Value2 = "super-liga"
$ws.Range("D124").Value2 = "2023-2024"
$ws.Range("E124").Value2 = 45262.79166666666
$ws.Range("F124").Value2 = "Vozdovac"
$ws.Range("G124").Value2 = 4
$ws.Range("H124").Value2 = "IMT Novi Beograd"
$ws.Range("I124").Value2 = 0
$ws.Range("J124").Value2 = 2.41
$ws.Range("K124").Value2 = "01/12/2023 07:12"
$ws.Range("L124").Value2 = 2.15
$ws.Range("M124").Value2 = "02/12/2023 18:27"
$ws.Range("N124").Value2 = 2.98
$ws.Range("O124").Value2 = "01/12/2023 07:12"
$ws.Range("P124").Value2 = 3.33
$ws.Range("Q124").Value2 = "02/12/2023 18:26"
$ws.Range("R124").Value2 = 2.8
$ws.Range("S124").Value2 = "01/12/2023 07:12"
$ws.Range("T124").Value2 = 3.27
$ws.Range("U124").Value2 = "02/12/2023 18:27"
$ws.Range("V124").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-imt-novi-beograd/8r7IAgHj/"
